$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "HK_G_acc_SD"

$ws.Range("A2").Value = 45.25745257452575
$ws.Range("A3").Value = 43.089430894308947
$ws.Range("A4").Value = 44.715447154471541
$ws.Range("A5").Value = 44.986449864498645
$ws.Range("A6").Value = 45.25745257452575
$ws.Range("A7").Value = 45.25745257452575
$ws.Range("A8").Value = 45.528455284552841
$ws.Range("A9").Value = 49.322493224932252
$ws.Range("A10").Value = 46.341463414634148
$ws.Range("A11").Value = 45.799457994579946
$ws.Range("A12").Value = 43.089430894308947
$ws.Range("A13").Value = 46.883468834688344
$ws.Range("A14").Value = 47.154471544715449
$ws.Range("A15").Value = 46.883468834688344
$ws.Range("A16").Value = 46.612466124661246
$ws.Range("A17").Value = 46.341463414634148
$ws.Range("A18").Value = 46.070460704607044
$ws.Range("A19").Value = 48.780487804878049
$ws.Range("A20").Value = 45.25745257452575
$ws.Range("A21").Value = 45.25745257452575
$ws.Range("A22").Value = 46.341463414634148
$ws.Range("A23").Value = 41.192411924119241
$ws.Range("A24").Value = 38.211382113821138
$ws.Range("A25").Value = 37.94037940379404
$ws.Range("A26").Value = 47.696476964769644
$ws.Range("A27").Value = 45.799457994579946
$ws.Range("A28").Value = 50.948509485094853
$ws.Range("A29").Value = 45.799457994579946
$ws.Range("A30").Value = 46.883468834688344
$ws.Range("A31").Value = 47.154471544715449
$ws.Range("A32").Value = 41.192411924119241
$ws.Range("A33").Value = 40.650406504065039
$ws.Range("A34").Value = 42.005420054200542
$ws.Range("A35").Value = 39.295392953929536
$ws.Range("A36").Value = 37.669376693766935
$ws.Range("A37").Value = 44.986449864498645
$ws.Range("A38").Value = 37.398373983739837
$ws.Range("A39").Value = 37.94037940379404
$ws.Range("A40").Value = 38.482384823848236
$ws.Range("A41").Value = 47.154471544715449
$ws.Range("A42").Value = 47.154471544715449
$ws.Range("A43").Value = 47.154471544715449
$ws.Range("A44").Value = 45.25745257452575
$ws.Range("A45").Value = 45.25745257452575
$ws.Range("A46").Value = 44.715447154471541
$ws.Range("A47").Value = 43.360433604336045
$ws.Range("A48").Value = 49.322493224932252
$ws.Range("A49").Value = 45.799457994579946
